$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force column D to Text format first so numeric-looking strings
# (e.g. "1.00", "215.57") are stored as text, matching the source
# workbook which keeps these as inline/shared strings, not numbers.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "26.620.99"
$ws.Range("D3").Value = "1.643.14"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "215.57"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "0.506"
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "19.20"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "1.667.49"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("E14").Value = "  +2.91%  "
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "65.76"
$ws.Range("E16").Value = "  +3.90%  "
$ws.Range("D17").Value = "26.653.56"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "0.0₃0749"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "217.57"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").Value = "9.52"
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("D24").Value = "2.16"
$ws.Range("E24").Value = "  +10.20%  "
$ws.Range("D25").Value = "146.07"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("E28").Value = "  +4.28%  "
$ws.Range("D29").Value = "15.80"
$ws.Range("E29").Value = "  +1.67%  "
$ws.Range("D30").Value = "0.0516"
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("D33").Value = "3.05"
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("D34").Value = "1.269.93"
$ws.Range("E34").Value = "  +4.67%  "
$ws.Range("E35").Value = "  +2.37%  "
$ws.Range("E36").Value = "  +5.34%  "
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("E38").Value = "  +5.71%  "
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").Value = "0.811"
$ws.Range("E41").Value = "  +2.47%  "
$ws.Range("D43").Value = "5.46"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("D44").Value = "1.782.30"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").Value = "92.98"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "59.43"
$ws.Range("E46").Value = "  +8.45%  "
$ws.Range("D47").Value = "1.60"
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("D49").Value = "7.81"
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("E51").Value = "  -0.36%  "

# Restore the default style on column D so we do not leave a stray
# number-format override on cells that did not have one originally.
$priceCol.Style = "Normal"
